$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 27; I = "sv"; J = "Statement-opinion" },
    @{ Row = 28; I = "aa"; J = "Agree/Accept" },
    @{ Row = 31; I = "aa"; J = "Agree/Accept" },
    @{ Row = 34; I = "qy"; J = "Yes-No-Question" },
    @{ Row = 36; I = "%"; J = "Uninterpretable" },
    @{ Row = 40; I = "%"; J = "Uninterpretable" },
    @{ Row = 42; I = "aa"; J = "Agree/Accept" },
    @{ Row = 44; I = "aa"; J = "Agree/Accept" },
    @{ Row = 80; I = "sv"; J = "Statement-opinion" },
    @{ Row = 90; I = "sv"; J = "Statement-opinion" },
    @{ Row = 109; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 121; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 124; I = "sv"; J = "Statement-opinion" },
    @{ Row = 126; I = "aa"; J = "Agree/Accept" },
    @{ Row = 129; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 131; I = "aa"; J = "Agree/Accept" },
    @{ Row = 142; I = "%"; J = "Uninterpretable" },
    @{ Row = 144; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 147; I = "aa"; J = "Agree/Accept" },
    @{ Row = 149; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 150; I = "aa"; J = "Agree/Accept" },
    @{ Row = 151; I = "aa"; J = "Agree/Accept" },
    @{ Row = 152; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 159; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 175; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 192; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 195; I = "%"; J = "Uninterpretable" },
    @{ Row = 196; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 197; I = "aa"; J = "Agree/Accept" },
    @{ Row = 198; I = "aa"; J = "Agree/Accept" },
    @{ Row = 210; I = "%"; J = "Uninterpretable" },
    @{ Row = 211; I = "%"; J = "Uninterpretable" },
    @{ Row = 213; I = "%"; J = "Uninterpretable" },
    @{ Row = 216; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 232; I = "%"; J = "Uninterpretable" },
    @{ Row = 234; I = "ba"; J = "Appreciation" },
    @{ Row = 236; I = "%"; J = "Uninterpretable" },
    @{ Row = 237; I = "%"; J = "Uninterpretable" },
    @{ Row = 251; I = "%"; J = "Uninterpretable" },
    @{ Row = 253; I = "aa"; J = "Agree/Accept" },
    @{ Row = 256; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 259; I = "%"; J = "Uninterpretable" },
    @{ Row = 260; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 262; I = "aa"; J = "Agree/Accept" },
    @{ Row = 265; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 267; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 281; I = "aa"; J = "Agree/Accept" },
    @{ Row = 282; I = "aa"; J = "Agree/Accept" },
    @{ Row = 289; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 296; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 301; I = "aa"; J = "Agree/Accept" },
    @{ Row = 310; I = "aa"; J = "Agree/Accept" },
    @{ Row = 311; I = "aa"; J = "Agree/Accept" },
    @{ Row = 332; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 344; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 356; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 363; I = "aa"; J = "Agree/Accept" },
    @{ Row = 371; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 379; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 385; I = "sv"; J = "Statement-opinion" },
    @{ Row = 387; I = "sv"; J = "Statement-opinion" },
    @{ Row = 398; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 399; I = "sv"; J = "Statement-opinion" },
    @{ Row = 405; I = "sv"; J = "Statement-opinion" },
    @{ Row = 433; I = "sv"; J = "Statement-opinion" },
    @{ Row = 437; I = "aa"; J = "Agree/Accept" },
    @{ Row = 459; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 475; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 478; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 488; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 494; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 495; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 499; I = "sv"; J = "Statement-opinion" },
    @{ Row = 502; I = "sv"; J = "Statement-opinion" },
    @{ Row = 509; I = "sv"; J = "Statement-opinion" },
    @{ Row = 534; I = "sv"; J = "Statement-opinion" },
    @{ Row = 542; I = "sv"; J = "Statement-opinion" },
    @{ Row = 545; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 546; I = "sd"; J = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
